# Swap the contents of rows 2-5 as a reversal: row2<->row5, row3<->row4.
# This matches the diff, where every column of each row pair was moved
# together (the rows were reordered, not individual fields edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51  # column AY

function Swap-Rows($rowA, $rowB) {
    for ($col = 1; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        # Only touch cells whose values actually differ. This avoids
        # pointlessly rewriting identical text (e.g. dates stored as
        # plain text like "2023-09-10"), which would otherwise risk
        # Excel reinterpreting/reformatting the text upon reassignment,
        # and it also leaves untouched/empty cells completely alone.
        if ($valA -ne $valB) {
            $cellA.Value = $valB
            $cellB.Value = $valA
        }
    }
}

Swap-Rows 2 5
Swap-Rows 3 4
